$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.17671077558947
$ws.Range("C2").Value = 0.1945268116595571
$ws.Range("D2").Value = 0.07786673192867966
$ws.Range("E2").Value = 0.06432310723595247
$ws.Range("G2").Value = 1.69463639135941
$ws.Range("H2").Value = 1.42132879432458
$ws.Range("L2").Value = 0.2500532008192664
$ws.Range("N2").Value = 1.590707076692588

$ws.Range("B3").Value = 2.0362873016403
$ws.Range("C3").Value = 0.1705232375492471
$ws.Range("D3").Value = 0.07070354049135119
$ws.Range("E3").Value = 0.06464240156176082
$ws.Range("G3").Value = 1.655235560946693
$ws.Range("H3").Value = 1.410139595656432
$ws.Range("L3").Value = 0.2403139995967791
$ws.Range("N3").Value = 1.609870562354004

$ws.Range("B4").Value = 1.951282567180442
$ws.Range("C4").Value = 0.1557365021969019
$ws.Range("D4").Value = 0.0663479394053752
$ws.Range("E4").Value = 0.0648589869073497
$ws.Range("G4").Value = 1.632282440568133
$ws.Range("H4").Value = 1.404128167953047
$ws.Range("L4").Value = 0.2344886492888492
$ws.Range("N4").Value = 1.622271959818818

$ws.Range("B5").Value = 1.916946316207884
$ws.Range("C5").Value = 0.1496982272453238
$ws.Range("D5").Value = 0.06458353665868799
$ws.Range("E5").Value = 0.06495241453690515
$ws.Range("G5").Value = 1.623238190215716
$ws.Range("H5").Value = 1.401893521516001
$ws.Range("L5").Value = 0.2321534564954106
$ws.Range("N5").Value = 1.62748483242386

$ws.Range("B6").Value = 1.911263107254115
$ws.Range("C6").Value = 0.1486948067374385
$ws.Range("D6").Value = 0.06429119033123243
$ws.Range("E6").Value = 0.06496824031080894
$ws.Range("G6").Value = 1.621755016558097
$ws.Range("H6").Value = 1.401535423435632
$ws.Range("L6").Value = 0.2317680307875492
$ws.Range("N6").Value = 1.628360028074319

$ws.Range("B7").Value = 1.950818268586318
$ws.Range("C7").Value = 0.1556551192549307
$ws.Range("D7").Value = 0.06632410160949576
$ws.Range("E7").Value = 0.06486022597840702
$ws.Range("G7").Value = 1.632159217160961
$ws.Range("H7").Value = 1.404097161156244
$ws.Range("L7").Value = 0.2344569996629247
$ws.Range("N7").Value = 1.622341618462357

$ws.Range("B8").Value = 2.128039230299692
$ws.Range("C8").Value = 0.1862601978244527
$ws.Range("D8").Value = 0.075387905792482
$ws.Range("E8").Value = 0.06442893918415482
$ws.Range("G8").Value = 1.680792416581681
$ws.Range("H8").Value = 1.417291972058052
$ws.Range("L8").Value = 0.2466629325043783
$ws.Range("N8").Value = 1.597182589857354

$ws.Range("B9").Value = 2.485321991590752
$ws.Range("C9").Value = 0.2459122685878867
$ws.Range("D9").Value = 0.09350927567453482
$ws.Range("E9").Value = 0.06374602347414893
$ws.Range("G9").Value = 1.78610367663353
$ws.Range("H9").Value = 1.450024704884356
$ws.Range("L9").Value = 0.2718348831269282
$ws.Range("N9").Value = 1.552903683195566

$ws.Range("B10").Value = 2.75394108017133
$ws.Range("C10").Value = 0.2895483635828668
$ws.Range("D10").Value = 0.10704869288098
$ws.Range("E10").Value = 0.063343420543994
$ws.Range("G10").Value = 1.869699360073952
$ws.Range("H10").Value = 1.478319725039398
$ws.Range("L10").Value = 0.2910982901874775
$ws.Range("N10").Value = 1.523481859438299

$ws.Range("B11").Value = 2.877513468385871
$ws.Range("C11").Value = 0.3093654393601355
$ws.Range("D11").Value = 0.1132602415792689
$ws.Range("E11").Value = 0.06318176834809641
$ws.Range("G11").Value = 1.909116742277035
$ws.Range("H11").Value = 1.492128595244822
$ws.Range("L11").Value = 0.3000324744861302
$ws.Range("N11").Value = 1.510778188732012

$ws.Range("B12").Value = 2.924507785612263
$ws.Range("C12").Value = 0.3168653801479877
$ws.Range("D12").Value = 0.1156201558775933
$ws.Range("E12").Value = 0.06312364387577496
$ws.Range("G12").Value = 1.924245632155589
$ws.Range("H12").Value = 1.497493570184872
$ws.Range("L12").Value = 0.3034404865342282
$ws.Range("N12").Value = 1.506066019612582

$ws.Range("B13").Value = 2.914377781220139
$ws.Range("C13").Value = 0.3152503206812014
$ws.Range("D13").Value = 0.1151115588423721
$ws.Range("E13").Value = 0.06313602462117629
$ws.Range("G13").Value = 1.92097831214943
$ws.Range("H13").Value = 1.496332067404666
$ws.Range("L13").Value = 0.3027054023156523
$ws.Range("N13").Value = 1.507076482178867

$ws.Range("B14").Value = 2.881375695106442
$ws.Range("C14").Value = 0.3099825488708916
$ws.Range("D14").Value = 0.1134542369561586
$ws.Range("E14").Value = 0.06317692449607648
$ws.Range("G14").Value = 1.910357334056954
$ws.Range("H14").Value = 1.492567246028528
$ws.Range("L14").Value = 0.3003123547891278
$ws.Range("N14").Value = 1.510388539567316

$ws.Range("B15").Value = 2.861187095418757
$ws.Range("C15").Value = 0.3067553318515479
$ws.Range("D15").Value = 0.1124400928541007
$ws.Range("E15").Value = 0.06320237919157989
$ws.Range("G15").Value = 1.903878112408876
$ws.Range("H15").Value = 1.490278908662106
$ws.Range("L15").Value = 0.2988497842575839
$ws.Range("N15").Value = 1.512430107063182

$ws.Range("B16").Value = 2.745893216921615
$ws.Range("C16").Value = 0.2882526385772906
$ws.Range("D16").Value = 0.1066438250668114
$ws.Range("E16").Value = 0.06335441731078717
$ws.Range("G16").Value = 1.867151519337625
$ws.Range("H16").Value = 1.477436235264065
$ws.Range("L16").Value = 0.2905178829413018
$ws.Range("N16").Value = 1.524325816691139

$ws.Range("B17").Value = 2.675518212644931
$ws.Range("C17").Value = 0.276893627118767
$ws.Range("D17").Value = 0.103101565214331
$ws.Range("E17").Value = 0.06345319205283673
$ws.Range("G17").Value = 1.844978616279434
$ws.Range("H17").Value = 1.469798497558742
$ws.Range("L17").Value = 0.2854505235095388
$ws.Range("N17").Value = 1.531798145887102

$ws.Range("B18").Value = 2.63516977031253
$ws.Range("C18").Value = 0.2703570734549317
$ws.Range("D18").Value = 0.1010690709973971
$ws.Range("E18").Value = 0.06351202783572152
$ws.Range("G18").Value = 1.832355957295988
$ws.Range("H18").Value = 1.465493633184565
$ws.Range("L18").Value = 0.2825520098534895
$ws.Range("N18").Value = 1.536160044972817

$ws.Range("B19").Value = 2.621530638159811
$ws.Range("C19").Value = 0.2681433549295775
$ws.Range("D19").Value = 0.1003817423229236
$ws.Range("E19").Value = 0.06353229613681499
$ws.Range("G19").Value = 1.828104495146391
$ws.Range("H19").Value = 1.4640511893829
$ws.Range("L19").Value = 0.2815733811492862
$ws.Range("N19").Value = 1.53764788565416

$ws.Range("B20").Value = 2.682996343687762
$ws.Range("C20").Value = 0.2781031363563216
$ws.Range("D20").Value = 0.1034781341810742
$ws.Range("E20").Value = 0.06344246793136143
$ws.Range("G20").Value = 1.8473254234313
$ws.Range("H20").Value = 1.470602415827443
$ws.Range("L20").Value = 0.2859882849877664
$ws.Range("N20").Value = 1.530996074382699

$ws.Range("B21").Value = 2.891063758034363
$ws.Range("C21").Value = 0.3115299359818948
$ws.Range("D21").Value = 0.113940821013415
$ws.Range("E21").Value = 0.06316482736570705
$ws.Range("G21").Value = 1.913471459236604
$ws.Range("H21").Value = 1.493669369635597
$ws.Range("L21").Value = 0.3010145752119087
$ws.Range("N21").Value = 1.509413031807313

$ws.Range("B22").Value = 3.028215465440326
$ws.Range("C22").Value = 0.3333510932866091
$ws.Range("D22").Value = 0.1208239759081806
$ws.Range("E22").Value = 0.06300138197666527
$ws.Range("G22").Value = 1.957882592581228
$ws.Range("H22").Value = 1.509537336150458
$ws.Range("L22").Value = 0.310979959151183
$ws.Range("N22").Value = 1.495881300458962

$ws.Range("B23").Value = 2.954907425796193
$ws.Range("C23").Value = 0.3217068953364617
$ws.Range("D23").Value = 0.1171461052228437
$ws.Range("E23").Value = 0.06308696837460914
$ws.Range("G23").Value = 1.9340706149151
$ws.Range("H23").Value = 1.500995443059196
$ws.Range("L23").Value = 0.3056479207580338
$ws.Range("N23").Value = 1.503050719991165

$ws.Range("B24").Value = 2.679615134806056
$ws.Range("C24").Value = 0.2775563362821174
$ws.Range("D24").Value = 0.103307874912332
$ws.Range("E24").Value = 0.06344730992482361
$ws.Range("G24").Value = 1.84626404327085
$ws.Range("H24").Value = 1.470238696050131
$ws.Range("L24").Value = 0.2857451170066838
$ws.Range("N24").Value = 1.53135848551112

$ws.Range("B25").Value = 2.387605238878336
$ws.Range("C25").Value = 0.2298105460947113
$ws.Range("D25").Value = 0.08856835540727559
$ws.Range("E25").Value = 0.06391335190003833
$ws.Range("G25").Value = 1.756533153576896
$ws.Range("H25").Value = 1.440428741667972
$ws.Range("L25").Value = 0.2648912295165786
$ws.Range("N25").Value = 1.564338118520453
